$d = $word.ActiveDocument

# 1. Fix capitalization: "openstack" -> "Openstack" in the "熟悉openstack" bullet item.
$d.Content.Find.Execute("熟悉openstack", $true, $false, $false, $false, $false,
                         $true, 1, $false, "熟悉Openstack", 2)

# 2. Move the "_GoBack" bookmark from the end of that paragraph to the empty
#    paragraph that immediately precedes the "技能证书" heading.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^\s*熟悉Openstack\r?$") {
        $target = $d.Paragraphs.Item($i + 1)
        $d.Bookmarks.Add("_GoBack", $target.Range)
        break
    }
}
